$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression: only C2 changes
$ws.Range("C2").Value = 0.7292314838023964

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 0.9964601137075352
$ws.Range("C3").Value = 0.9964937946425758
$ws.Range("D3").Value = 0.9785968291125831

# Row 4 - label change from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9960721656778091
$ws.Range("C4").Value = 0.9962328787497309
$ws.Range("D4").Value = 0.9580386068868422

# Row 5 - label change from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9976341310443648
$ws.Range("C5").Value = 0.9974299833969115
$ws.Range("D5").Value = 0.9969779540689164
